$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CompStat_1")

# Update the volume/date header text (rich-text shared strings)
$ws.Range("A8").Value = "Volume 31   Number  23"
$ws.Range("C9").Value = "Report Covering the Week  6/3/2024  Through  6/9/2024"

# Cells whose number-format class changes: copy formats from a stable donor cell first
$ws.Range("F31").Copy()
$ws.Range("C14").PasteSpecial(-4122)
$ws.Range("F31").Copy()
$ws.Range("F14").PasteSpecial(-4122)
$ws.Range("F31").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("F31").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("H31").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("C31").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("F31").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("F31").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("F31").Copy()
$ws.Range("C29").PasteSpecial(-4122)
$ws.Range("F31").Copy()
$ws.Range("C30").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Set cell values
$ws.Range("C14").Value = 1
$ws.Range("F14").Value = 1
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 2
$ws.Range("K14").Value = -33.333333333333
$ws.Range("L14").Value = 100
$ws.Range("M14").Value = 0
$ws.Range("N14").Value = -50
$ws.Range("C15").Value = 2
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 100
$ws.Range("F15").Value = 3
$ws.Range("G15").Value = 3
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 13
$ws.Range("J15").Value = 8
$ws.Range("K15").Value = 62.5
$ws.Range("L15").Value = 85.714285714285
$ws.Range("M15").Value = 116.666666666667
$ws.Range("N15").Value = -7.142857142857
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = -66.666666666666
$ws.Range("G16").Value = 12
$ws.Range("H16").Value = 16.666666666666
$ws.Range("J16").Value = 67
$ws.Range("K16").Value = 17.910447761194
$ws.Range("L16").Value = -12.222222222222
$ws.Range("M16").Value = -32.478632478632
$ws.Range("N16").Value = -84.136546184739
$ws.Range("C17").Value = 10
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = 150
$ws.Range("F17").Value = 32
$ws.Range("H17").Value = 60
$ws.Range("I17").Value = 153
$ws.Range("J17").Value = 159
$ws.Range("K17").Value = -3.773584905660
$ws.Range("L17").Value = -8.383233532934
$ws.Range("M17").Value = 62.765957446808
$ws.Range("N17").Value = 2
$ws.Range("C18").Value = "'0"
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = -100
$ws.Range("F18").Value = 3
$ws.Range("G18").Value = 13
$ws.Range("H18").Value = -76.923076923076
$ws.Range("J18").Value = 69
$ws.Range("K18").Value = -42.028985507246
$ws.Range("L18").Value = -24.528301886792
$ws.Range("M18").Value = -69.465648854961
$ws.Range("N18").Value = -93.265993265993
$ws.Range("C19").Value = 3
$ws.Range("D19").Value = 11
$ws.Range("E19").Value = -72.727272727272
$ws.Range("G19").Value = 35
$ws.Range("H19").Value = -34.285714285714
$ws.Range("I19").Value = 142
$ws.Range("J19").Value = 141
$ws.Range("K19").Value = 0.709219858156
$ws.Range("L19").Value = -10.691823899371
$ws.Range("M19").Value = 1.428571428571
$ws.Range("N19").Value = -40.585774058577
$ws.Range("D20").Value = 7
$ws.Range("E20").Value = -57.142857142857
$ws.Range("F20").Value = 19
$ws.Range("G20").Value = 26
$ws.Range("H20").Value = -26.923076923076
$ws.Range("I20").Value = 114
$ws.Range("J20").Value = 119
$ws.Range("K20").Value = -4.201680672268
$ws.Range("L20").Value = 5.555555555555
$ws.Range("M20").Value = 8.571428571428
$ws.Range("N20").Value = -92.738853503184
$ws.Range("C21").Value = 20
$ws.Range("D21").Value = 30
$ws.Range("E21").Value = -33.333333333333
$ws.Range("G21").Value = 110
$ws.Range("H21").Value = -13.636363636363
$ws.Range("I21").Value = 543
$ws.Range("J21").Value = 566
$ws.Range("K21").Value = -4.063604240282
$ws.Range("L21").Value = -7.179487179487
$ws.Range("M21").Value = -8.739495798319
$ws.Range("N21").Value = -82.306940371456
$ws.Range("M22").Value = -62.5
$ws.Range("D24").Value = 34
$ws.Range("E24").Value = -47.058823529411
$ws.Range("F24").Value = 102
$ws.Range("G24").Value = 128
$ws.Range("H24").Value = -20.3125
$ws.Range("I24").Value = 623
$ws.Range("J24").Value = 661
$ws.Range("K24").Value = -5.748865355521
$ws.Range("L24").Value = -2.044025157232
$ws.Range("M24").Value = 71.153846153846
$ws.Range("D25").Value = 12
$ws.Range("E25").Value = -41.666666666666
$ws.Range("F25").Value = 45
$ws.Range("G25").Value = 54
$ws.Range("H25").Value = -16.666666666666
$ws.Range("I25").Value = 288
$ws.Range("J25").Value = 216
$ws.Range("K25").Value = 33.333333333333
$ws.Range("L25").Value = 15.2
$ws.Range("C26").Value = 8
$ws.Range("D26").Value = 10
$ws.Range("E26").Value = -20
$ws.Range("F26").Value = 37
$ws.Range("G26").Value = 54
$ws.Range("H26").Value = -31.481481481481
$ws.Range("I26").Value = 265
$ws.Range("J26").Value = 241
$ws.Range("K26").Value = 9.958506224066
$ws.Range("L26").Value = 25
$ws.Range("M26").Value = -0.375939849624
$ws.Range("C27").Value = 2
$ws.Range("E27").Value = 100
$ws.Range("F27").Value = 5
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = -16.666666666666
$ws.Range("I27").Value = 17
$ws.Range("J27").Value = 18
$ws.Range("K27").Value = -5.555555555555
$ws.Range("L27").Value = 6.25
$ws.Range("C28").Value = 2
$ws.Range("D28").Value = 5
$ws.Range("E28").Value = -60
$ws.Range("F28").Value = 5
$ws.Range("G28").Value = 13
$ws.Range("H28").Value = -61.538461538461
$ws.Range("I28").Value = 18
$ws.Range("J28").Value = 33
$ws.Range("K28").Value = -45.454545454545
$ws.Range("L28").Value = -14.285714285714
$ws.Range("C29").Value = 2
$ws.Range("F29").Value = 3
$ws.Range("I29").Value = 9
$ws.Range("K29").Value = 125
$ws.Range("L29").Value = 12.5
$ws.Range("M29").Value = 200
$ws.Range("N29").Value = -40
$ws.Range("C30").Value = 1
$ws.Range("F30").Value = 2
$ws.Range("I30").Value = 5
$ws.Range("K30").Value = 25
$ws.Range("L30").Value = -16.666666666666
$ws.Range("M30").Value = 66.666666666666
$ws.Range("N30").Value = -54.545454545454
